$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the two new performance test rows' values to existing row 21, and create new rows 22-26 ---
# Use existing formatted rows as templates so number/fill styles (s=2 / s=1 / s=3) are copied correctly.

# Row 21 already has A21 filled in; just add B21/C21 values (copy B20:C20 formatting first).
$ws.Range("B20:C20").Copy($ws.Range("B21:C21"))
$ws.Range("B21").Value = 1500
$ws.Range("C21").Value = 5000

# Row 22 : RegistrarBanner
$ws.Range("A20:C20").Copy($ws.Range("A22:C22"))
$ws.Range("A22").Value = "RegistrarBanner"
$ws.Range("B22").Value = 300
$ws.Range("C22").Value = 2500
$ws.Range("E22").Value = "Profiling"

# Row 23 : EliminarBanners
$ws.Range("A20:C20").Copy($ws.Range("A23:C23"))
$ws.Range("A23").Value = "EliminarBanners"
$ws.Range("B23").Value = 300
$ws.Range("C23").Value = 2500
$ws.Range("E23").Value = "Profiling"

# Row 24 : ActualizarDatosDeUsuario
$ws.Range("A20:C20").Copy($ws.Range("A24:C24"))
$ws.Range("A24").Value = "ActualizarDatosDeUsuario"
$ws.Range("B24").Value = 3000
$ws.Range("C24").Value = 7000

# Row 25 : SolicitarRegistroMascota
$ws.Range("A20:C20").Copy($ws.Range("A25:C25"))
$ws.Range("A25").Value = "SolicitarRegistroMascota"
$ws.Range("B25").Value = 50
$ws.Range("C25").Value = 300
$ws.Range("E25").Value = "Profiling"

# Row 26 : AceptarRechazarSolicitudMascota
$ws.Range("A20:C20").Copy($ws.Range("A26:C26"))
$ws.Range("A26").Value = "AceptarRechazarSolicitudMascota"
$ws.Range("B26").Value = 10
$ws.Range("C26").Value = 500
$ws.Range("E26").Value = "Profiling"

# --- Highlight (yellow fill) the "Profiling" cells that were marked in the diff ---
$ws.Range("E18").Interior.Color = 65535
$ws.Range("E19").Interior.Color = 65535
$ws.Range("E25").Interior.Color = 65535
$ws.Range("E26").Interior.Color = 65535

# --- Remove the two obsolete "Profiling" markers ---
$ws.Range("E14").ClearContents()
$ws.Range("E15").ClearContents()

# --- Widen column A to fit the longer new user-story names ---
$ws.Columns.Item(1).ColumnWidth = 32.66

# --- Update the selection to match the final cursor position recorded in the workbook ---
$ws.Range("G22").Select()
